$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (cultivo-descripcion), H (grupo-cultivo-descripcion) and J (secanoregadio)
# move from curated "dimension" metadata (dim / skos:Concept / mapping file) to
# "measure" metadata (medida / xsd:int, no mapping file needed anymore).
foreach ($col in @("D", "H", "J")) {
    $ws.Range($col + "2").Value = "iaest-measure:" + $ws.Range($col + "1").Value2
    $ws.Range($col + "3").Value = "medida"
    $ws.Range($col + "4").Value = "xsd:int"
    $ws.Range($col + "5").Clear()
}

# Column K (municipio-nombre) stops being a plain measure and instead becomes a
# referenced area dimension, like columns L (provincia-nombre) and M (comarca-nombre).
$ws.Range("K2").Value = "sdmx-dimension:refArea"
$ws.Range("K3").Value = "dim"
$ws.Range("K4").Value = "URI-Municipio"
